$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value = "Authorship Resource"
$ws.Range("I2").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I3").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I4").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I5").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("J15").Select()
